$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventory for check stock")

# Update the Location No. for the row whose Item Code is 3000041389
# (เครื่องถอนขนถนอมผิว EMJOI LIGHT) from "01A012" to "01A000"
$ws.Range("F6").Value = "01A000"

# Move the active selection to D7, as recorded in the saved sheet view
$ws.Range("D7").Select()
